$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: new Price (D) and Volume(1h) (E) text values.
# NumberFormat "@" forces text interpretation so numeric-looking strings
# (e.g. "1.29") are not auto-coerced to numbers by Excel; ClearFormats()
# afterwards restores the default (unstyled) cell appearance.
$updates = @{
    'D2' = '72.241.76'
    'E2' = '  +1.35%  '
    'D3' = '2.671.13'
    'E3' = '  +1.52%  '
    'D5' = '599.06'
    'E5' = '  -1.34%  '
    'D6' = '176.07'
    'E6' = '  -2.75%  '
    'E7' = '  +0.01%  '
    'E8' = '  -0.61%  '
    'D9' = '2.670.30'
    'E9' = '  +1.54%  '
    'D10' = '0.169'
    'E10' = '  +0.93%  '
    'E11' = '  +2.11%  '
    'E12' = '  +1.78%  '
    'E13' = '  +0.31%  '
    'D14' = '3.164.39'
    'E14' = '  +1.61%  '
    'E15' = '  -2.66%  '
    'D16' = '72.128.67'
    'E16' = '  +1.35%  '
    'D17' = '26.25'
    'E17' = '  -1.73%  '
    'D18' = '2.680.36'
    'E18' = '  +0.91%  '
    'D19' = '12.02'
    'E19' = '  +4.22%  '
    'D20' = '7.97'
    'E20' = '  +0.63%  '
    'D21' = '370.01'
    'E21' = '  -3.16%  '
    'D22' = '4.15'
    'E22' = '  +0.26%  '
    'E23' = '  +4.38%  '
    'D24' = '71.76'
    'E24' = '  -1.07%  '
    'E25' = '  +0.00%  '
    'D26' = '4.33'
    'E26' = '  -3.40%  '
    'D27' = '9.80'
    'E27' = '  +1.29%  '
    'D28' = '2.808.70'
    'E28' = '  +1.54%  '
    'E29' = '  +0.16%  '
    'D30' = '0.0₃0937'
    'E30' = '  -3.28%  '
    'D31' = '8.04'
    'E31' = '  -0.31%  '
    'D32' = '509.78'
    'E32' = '  -6.57%  '
    'D33' = '1.29'
    'E33' = '  -2.41%  '
    'E34' = '  -1.46%  '
    'E35' = '  +0.02%  '
    'D36' = '164.79'
    'E36' = '  -0.79%  '
    'E37' = '  +1.65%  '
    'E38' = '  +0.29%  '
    'D39' = '1.37'
    'E39' = '  -0.93%  '
    'E40' = '  -3.76%  '
    'D41' = '0.105'
    'E41' = '  -10.60%  '
    'E42' = '  -0.05%  '
    'D43' = '5.00'
    'E43' = '  -1.12%  '
    'D44' = '2.56'
    'E44' = '  -3.16%  '
    'E45' = '  -0.03%  '
    'D46' = '39.24'
    'E46' = '  -1.81%  '
    'D47' = '153.34'
    'E47' = '  -0.81%  '
    'E48' = '  +2.00%  '
    'D49' = '0.548'
    'E49' = '  +2.47%  '
    'E50' = '  +1.24%  '
    'E51' = '  +1.31%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
